$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-150) holds the "Förändrad" (changed) date, stored as the
# Excel date serial number 45175 (2023-09-06). Update it to 45177
# (2023-09-08) for every data row.
$ws.Range("C2:C150").Value = 45177
